$wb = $excel.ActiveWorkbook

# --- Sheet "phpmyadmin" (sheet1) ---
$ws1 = $wb.Worksheets.Item("phpmyadmin")

# Row 15 "total" row gets taller (BFIAF -> FFIAF generalization bumped its wrapped header text)
$ws1.Rows.Item(15).RowHeight = 25

# Selection moves from R4 to K2 (and this sheet is no longer the tab shown when reopened)
$ws1.Range("K2").Select()

# --- Sheet "drupal" (sheet2) ---
$ws2 = $wb.Worksheets.Item("drupal")

# New column O was generated (VFIAF Images) - widen it to fit its contents
$ws2.Columns.Item(15).ColumnWidth = 8.5

# Row 16 "total" row gets taller too
$ws2.Rows.Item(16).RowHeight = 25

# Selection moves to A15
$ws2.Range("A15").Select()

# --- Sheet "moodle" (sheet3) ---
$ws3 = $wb.Worksheets.Item("moodle")

# Generalized BFIAF to FFIAF: the totals row now sums the whole column (from row 2)
# instead of just the BFIAF-image subset (rows 23:35).
$ws3.Range("C36:AR36").Formula = "=SUM(C2:C35)"
$ws3.Range("B36").Formula = "=SUM(B2:B35)"

# Total row is taller to match the other regenerated sheets
$ws3.Rows.Item(36).RowHeight = 37

# moodle becomes the active/selected tab, with B36 selected
$ws3.Range("B36").Select()
$ws3.Activate()
